{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"Docente(s) Respons\u00e1vel(eis)\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find 'Docente(s) Respons\u00e1vel(eis)' paragraph\");\n}\n\nconst newPara = target.insertParagraph(\"7455355 - Robson da Silva Rocha\", \"After\");\nnewPara.style = \"ListBullet\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Docente(s) Respons\u00e1vel(eis)*\") {\n        $newRange = $p.Range.InsertParagraphAfter()\n        $newPara = $p.Next()\n        $newPara.Range.Text = \"7455355 - Robson da Silva Rocha\"\n        $newPara.Style = \"ListBullet\"\n        break\n    }\n}\n"}
